# Update the header_table sheet's style-guide values:
#  - font_family cell for TITLE row gets quoted family + sans-serif fallback
#  - font_size cells (px unit) for TITLE/PRODUCT/NUMBER rows
#  - move active selection back to B1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("header_table")

$ws.Range("C1").Value = "font_size: 14px"
$ws.Range("C2").Value = "font_size: 14px"
$ws.Range("C3").Value = "font_size: 14px"
$ws.Range("B1").Value = 'font_family: "Open Sans", sans-serif'

$ws.Activate()
$ws.Range("B1").Select()
